$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two answer names (shared strings "Bob" -> "Name1", "Kari" -> "Name2")
$ws.Range("B2").Value = "Name1"
$ws.Range("B3").Value = "Name2"

# Update the saved selection to B6
$ws.Range("B6").Select()
